$wb = $excel.ActiveWorkbook

$oldBuild = "January 30 2026 16.19.47 EST"
$newBuild = "February 02 2026 12.49.33 EST"

$oldVersion = "mines - January 30 (built on $oldBuild)"
$newVersion = "mines - January 30 (built on $newBuild)"

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: $newVersion"

$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Moranbah North Coal Mine, Australia, M0074, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$usedRange = $wsData.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsData.Cells.Item($r, 19)  # column S
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
